$wb = $excel.ActiveWorkbook

$oldFileName = "1d8dfc77-f9a0-4b80-aa22-6aa379968115.md"
$newFileName = "7eafc53e-b21e-443d-b489-0c20752e608d.md"

$newSourceFileName = "7eafc53e-b21e-443d-b489-0c20752e608d.md"
$newPathAndName = "e2e\7eafc53e-b21e-443d-b489-0c20752e608d.md"

$newLatestHoDate = "2016-08-21 13:01:54"

$newZhCnXlf = "7eafc53e-b21e-443d-b489-0c20752e608d.4bf2d6cb022035ed391ddf1fc2651833277492cc.zh-cn.xlf"
$newZhCnHandoffDate = "2016-08-21 13:01:50"

$newDeDeXlf = "7eafc53e-b21e-443d-b489-0c20752e608d.4bf2d6cb022035ed391ddf1fc2651833277492cc.de-de.xlf"
$newDeDeHandoffDate = "2016-08-21 13:01:54"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91e061b16cd2d26436a140bc0472103318b14990/e2e/1d8dfc77-f9a0-4b80-aa22-6aa379968115.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newSourceFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = $newLatestHoDate

$rB2 = $wsOverview.Range("B2")
$rB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rB2, $hyperlinkAddress, "", "", $newPathAndName)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newSourceFileName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate

$rZhA2 = $wsZhCn.Range("A2")
$rZhA2.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($rZhA2, $hyperlinkAddress, "", "", $newSourceFileName)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newSourceFileName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newDeDeHandoffDate

$rDeA2 = $wsDeDe.Range("A2")
$rDeA2.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($rDeA2, $hyperlinkAddress, "", "", $newSourceFileName)
